$wb = $excel.ActiveWorkbook

# --- Sheet: Login ---
$ws = $wb.Worksheets.Item("Login")
$ws.Range("E5").Value = "successfully"

# --- Sheet: Register ---
$ws = $wb.Worksheets.Item("Register")
$ws.Range("E6").Value = "+6281252930362"
$ws.Range("E7").Value = "+6281252930362"
$ws.Range("E8").Value = "+6281252930362"
$ws.Range("E9").Value = "+6281252930362"
$ws.Range("D10").Value = "testregisterbackend23@gmail.com"
$ws.Range("D11").Value = "testregisterbackend23@gmail.com"
$ws.Range("D12").Value = "testregisterbackend23@gmail.com"
$ws.Range("D13").Value = "testregisterbackend23@gmail.com"
$ws.Range("D14").Value = "testregisterbackend23@gmail.com"
$ws.Range("E14").Value = "+6281252930362"
$ws.Range("D15").Value = "testregisterbackend23@gmail.com"
$ws.Range("E15").Value = "+6281252930362"
$ws.Range("D16").Value = "testregisterbackend23@gmail.com"
$ws.Range("E16").Value = "+6281252930362"
$ws.Range("D17").Value = "testregisterbackend23@gmail.com"
$ws.Range("E17").Value = "+6281252930362"
$ws.Range("E22").Value = "+6281252930362"
$ws.Range("E23").Value = "+6281252930362"
$ws.Range("E24").Value = "+6281252930362"
$ws.Range("E25").Value = "+6281252930362"
$ws.Range("D26").Value = "testregisterbackend23@gmail.com"
$ws.Range("D27").Value = "testregisterbackend23@gmail.com"
$ws.Range("D28").Value = "testregisterbackend23@gmail.com"
$ws.Range("D29").Value = "testregisterbackend23@gmail.com"
$ws.Range("D30").Value = "testregisterbackend23@gmail.com"
$ws.Range("E30").Value = "+6281252930362"
$ws.Range("D31").Value = "testregisterbackend23@gmail.com"
$ws.Range("E31").Value = "+6281252930362"
$ws.Range("D32").Value = "testregisterbackend23@gmail.com"
$ws.Range("E32").Value = "+6281252930362"
$ws.Range("D33").Value = "testregisterbackend23@gmail.com"
$ws.Range("E33").Value = "+6281252930362"
$ws.Range("E38").Value = "+6281252930362"
$ws.Range("E39").Value = "+6281252930362"
$ws.Range("E40").Value = "+6281252930362"
$ws.Range("E41").Value = "+6281252930362"
$ws.Range("D42").Value = "testregisterbackend23@gmail.com"
$ws.Range("D43").Value = "testregisterbackend23@gmail.com"
$ws.Range("D44").Value = "testregisterbackend23@gmail.com"
$ws.Range("D45").Value = "testregisterbackend23@gmail.com"
$ws.Range("D46").Value = "testregisterbackend23@gmail.com"
$ws.Range("E46").Value = "+6281252930362"
$ws.Range("D47").Value = "testregisterbackend23@gmail.com"
$ws.Range("E47").Value = "+6281252930362"
$ws.Range("D48").Value = "testregisterbackend23@gmail.com"
$ws.Range("E48").Value = "+6281252930362"
$ws.Range("D49").Value = "testregisterbackend23@gmail.com"
$ws.Range("E49").Value = "+6281252930362"
$ws.Range("E54").Value = "+6281252930362"
$ws.Range("E55").Value = "+6281252930362"
$ws.Range("E56").Value = "+6281252930362"
$ws.Range("E57").Value = "+6281252930362"
$ws.Range("D58").Value = "testregisterbackend23@gmail.com"
$ws.Range("D59").Value = "testregisterbackend23@gmail.com"
$ws.Range("D60").Value = "testregisterbackend23@gmail.com"
$ws.Range("D61").Value = "testregisterbackend23@gmail.com"
$ws.Range("D62").Value = "testregisterbackend23@gmail.com"
$ws.Range("E62").Value = "+6281252930362"
$ws.Range("D63").Value = "testregisterbackend23@gmail.com"
$ws.Range("E63").Value = "+6281252930362"
$ws.Range("D64").Value = "testregisterbackend23@gmail.com"
$ws.Range("E64").Value = "+6281252930362"
$ws.Range("D65").Value = "testregisterbackend23@gmail.com"
$ws.Range("E65").Value = "+6281252930362"
$ws.Range("I65").Value = "+6281252930362 ----"
$ws.Range("O65").Value = "signup is successfully"
$ws.Range("D66").Value = "testregisterbackend23@gmail.com"
$ws.Range("E66").Value = "+6281252930362"
$ws.Range("D67").Value = "testregisterbackend23@gmail.com"
$ws.Range("E67").Value = "+6281252930362"
$ws.Range("D68").Value = "testregisterbackend23@gmail.com"
$ws.Range("E68").Value = "+6281252930362"
$ws.Range("D69").Value = "testregisterbackend23@gmail.com"
$ws.Range("E69").Value = "+6281252930362"
$ws.Range("D70").Value = "testregisterbackend23@gmail.com"
$ws.Range("E70").Value = "+6281252930362"
$ws.Range("I70").Value = ""
$ws.Range("D71").Value = "testregisterbackend23@gmail.com"
$ws.Range("E71").Value = "+6281252930362"
$ws.Range("D72").Value = "testregister.gmail.com"
$ws.Range("E72").Value = "+6281252930362"
$ws.Range("P72").Value = "SUCCESS"
$ws.Range("D73").Value = "testregisterbackend23@gmail.com"
$ws.Range("D74").Value = "testregisterbackend23@gmail.com"
$ws.Range("D75").Value = "testregisterbackend23@gmail.com"
$ws.Range("D76").Value = "testregisterbackend23@gmail.com"
$ws.Range("D77").Value = "testregisterbackend23@gmail.com"
$ws.Range("D78").Value = "testregisterbackend23@gmail.com"
$ws.Range("E78").Value = "+6281252930362"
$ws.Range("D79").Value = "testregisterbackend23@gmail.com"
$ws.Range("E79").Value = "+6281252930362"
$ws.Range("D80").Value = "testregisterbackend23@gmail.com"
$ws.Range("E80").Value = "+6281252930362"
$ws.Range("D81").Value = "testregisterbackend23@gmail.com"
$ws.Range("E81").Value = "+6281252930362"
$ws.Range("D82").Value = "testregisterbackend23@gmail.com"
$ws.Range("E82").Value = "+6281252930362"
$ws.Range("D83").Value = "testregisterbackend23@gmail.com"
$ws.Range("E83").Value = "+6281252930362"
$ws.Range("D84").Value = "testregisterbackend23@gmail.com"
$ws.Range("E84").Value = "+6281252930362"
$ws.Range("D85").Value = "testregisterbackend23@gmail.com"
$ws.Range("E85").Value = "+6281252930362"
$ws.Range("D86").Value = "testregisterbackend23@gmail.com"
$ws.Range("E86").Value = "+6281252930362"
$ws.Range("D87").Value = "testregisterbackend23@gmail.com"
$ws.Range("E87").Value = "+6281252930362"
$ws.Range("D88").Value = "testregisterbackend23@gmail.com"
$ws.Range("E88").Value = "+6281252930362"
$ws.Range("D89").Value = "testregisterbackend23@gmail.com"
$ws.Range("E89").Value = "+6281252930362"
$ws.Range("D90").Value = "testregisterbackend23@gmail.com"
$ws.Range("E90").Value = "+6281252930362"
$ws.Range("I90").Value = "+6281252930362 ----"
$ws.Range("D91").Value = "testregisterbackend23@gmail.com"
$ws.Range("E91").Value = "+6281252930362"
$ws.Range("I91").Value = "+6281252930362 ----"
$ws.Range("D92").Value = "testregisterbackend23@gmail.com"
$ws.Range("E92").Value = "+6281252930362"
$ws.Range("I92").Value = "+6281252930362 ----"
$ws.Range("D93").Value = "testregisterbackend23@gmail.com"
$ws.Range("E93").Value = "+6281252930362"
$ws.Range("I93").Value = "+6281252930362 ----"
$ws.Range("O93").Value = "signup is successfully"
$ws.Range("D94").Value = "testregisterbackend23@gmail.com"
$ws.Range("E94").Value = "+6281252930362"
$ws.Range("I94").Value = "+6281252930362 ----"
$ws.Range("D95").Value = "testregisterbackend23@gmail.com"
$ws.Range("E95").Value = "+6281252930362"
$ws.Range("I95").Value = "+6281252930362 ----"
$ws.Range("D96").Value = "testregisterbackend23@gmail.com"
$ws.Range("E96").Value = "+6281252930362"
$ws.Range("I96").Value = "+6281252930362 ----"
$ws.Range("D97").Value = "testregisterbackend23@gmail.com"
$ws.Range("E97").Value = "+6281252930362"
$ws.Range("I97").Value = "+6281252930362 ----"
$ws.Range("D98").Value = "testregisterbackend23@gmail.com"
$ws.Range("E98").Value = "+6281252930362"
$ws.Range("I98").Value = "+6281252930362 ----"
$ws.Range("P98").Value = "SUCCESS"
$ws.Range("D99").Value = "testregisterbackend23@gmail.com"
$ws.Range("E99").Value = "+6281252930362"
$ws.Range("I99").Value = "+6281252930362 ----"
$ws.Range("O99").Value = "signup is successfully"
$ws.Range("P99").Value = "SUCCESS"
$ws.Range("D100").Value = "testregisterbackend23@gmail.com"
$ws.Range("E100").Value = "+6281252930362"
$ws.Range("I100").Value = "+6281252930362 ----"
$ws.Range("O100").Value = "No message available"
$ws.Range("P100").Value = "SUCCESS"

# --- Sheet: Paket Data Purchase ---
$ws = $wb.Worksheets.Item("Paket Data Purchase")
$ws.Range("N53").Value = "Anda Tidak memiliki transaksi"
$ws.Range("N54").Value = "Anda Tidak memiliki transaksi"
$ws.Range("N55").Value = "Anda Tidak memiliki transaksi"
$ws.Range("Q62").Value = "Transacstion successfully"
$ws.Range("N65").Value = "Anda Tidak memiliki transaksi"
$ws.Range("N66").Value = "Anda Tidak memiliki transaksi"
$ws.Range("N67").Value = "Anda Tidak memiliki transaksi"
